$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns G:K (numeric-looking text) to remain text, matching the
# source data which stores these values as text, not numbers.
$ws.Range("G11:K19").NumberFormat = "@"

# Row 11
$ws.Range("A11").Value = ' Dubai (DSC)'
$ws.Range("B11").Value = ' October 14 2020'
$ws.Range("C11").Value = 'Capitals won by 13 runs'
$ws.Range("D11").Value = 'Delhi Capitals'
$ws.Range("E11").Value = 'Rajasthan Royals'
$ws.Range("F11").Value = 'Axar Patel '
$ws.Range("G11").Value = '7'
$ws.Range("H11").Value = '4'
$ws.Range("I11").Value = '1'
$ws.Range("J11").Value = '0'
$ws.Range("K11").Value = '175.00'

# Row 12
$ws.Range("A12").Value = ' Dubai (DSC)'
$ws.Range("B12").Value = ' November 05 2020'
$ws.Range("C12").Value = 'Mumbai won by 57 runs'
$ws.Range("D12").Value = 'Delhi Capitals'
$ws.Range("E12").Value = 'Mumbai Indians'
$ws.Range("F12").Value = 'Axar Patel '
$ws.Range("G12").Value = '42'
$ws.Range("H12").Value = '33'
$ws.Range("I12").Value = '2'
$ws.Range("J12").Value = '3'
$ws.Range("K12").Value = '127.27'

# Row 13
$ws.Range("A13").Value = ' Dubai (DSC)'
$ws.Range("B13").Value = ' October 27 2020'
$ws.Range("C13").Value = 'Sunrisers won by 88 runs'
$ws.Range("D13").Value = 'Delhi Capitals'
$ws.Range("E13").Value = 'Sunrisers Hyderabad'
$ws.Range("F13").Value = 'Axar Patel '
$ws.Range("G13").Value = '1'
$ws.Range("H13").Value = '4'
$ws.Range("I13").Value = '0'
$ws.Range("J13").Value = '0'
$ws.Range("K13").Value = '25.00'

# Row 14
$ws.Range("A14").Value = ' Abu Dhabi'
$ws.Range("B14").Value = ' September 29 2020'
$ws.Range("C14").Value = 'Sunrisers won by 15 runs'
$ws.Range("D14").Value = 'Delhi Capitals'
$ws.Range("E14").Value = 'Sunrisers Hyderabad'
$ws.Range("F14").Value = 'Axar Patel '
$ws.Range("G14").Value = '5'
$ws.Range("H14").Value = '6'
$ws.Range("I14").Value = '0'
$ws.Range("J14").Value = '0'
$ws.Range("K14").Value = '83.33'

# Row 15
$ws.Range("A15").Value = ' Abu Dhabi'
$ws.Range("B15").Value = ' October 24 2020'
$ws.Range("C15").Value = 'KKR won by 59 runs'
$ws.Range("D15").Value = 'Delhi Capitals'
$ws.Range("E15").Value = 'Kolkata Knight Riders'
$ws.Range("F15").Value = 'Axar Patel '
$ws.Range("G15").Value = '9'
$ws.Range("H15").Value = '7'
$ws.Range("I15").Value = '0'
$ws.Range("J15").Value = '1'
$ws.Range("K15").Value = '128.57'

# Row 16
$ws.Range("A16").Value = ' Sharjah'
$ws.Range("B16").Value = ' October 17 2020'
$ws.Range("C16").Value = 'Capitals won by 5 wickets (with 1 ball remaining)'
$ws.Range("D16").Value = 'Delhi Capitals'
$ws.Range("E16").Value = 'Chennai Super Kings'
$ws.Range("F16").Value = 'Axar Patel '
$ws.Range("G16").Value = '21'
$ws.Range("H16").Value = '5'
$ws.Range("I16").Value = '0'
$ws.Range("J16").Value = '3'
$ws.Range("K16").Value = '420.00'

# Row 17
$ws.Range("A17").Value = ' Dubai (DSC)'
$ws.Range("B17").Value = ' November 10 2020'
$ws.Range("C17").Value = 'Mumbai won by 5 wickets (with 8 balls remaining)'
$ws.Range("D17").Value = 'Delhi Capitals'
$ws.Range("E17").Value = 'Mumbai Indians'
$ws.Range("F17").Value = 'Axar Patel '
$ws.Range("G17").Value = '9'
$ws.Range("H17").Value = '9'
$ws.Range("I17").Value = '1'
$ws.Range("J17").Value = '0'
$ws.Range("K17").Value = '100.00'

# Row 18
$ws.Range("A18").Value = ' Dubai (DSC)'
$ws.Range("B18").Value = ' September 20 2020'
$ws.Range("C18").Value = 'Match tied (Capitals won the one-over eliminator)'
$ws.Range("D18").Value = 'Delhi Capitals'
$ws.Range("E18").Value = 'Kings XI Punjab'
$ws.Range("F18").Value = 'Axar Patel '
$ws.Range("G18").Value = '6'
$ws.Range("H18").Value = '9'
$ws.Range("I18").Value = '0'
$ws.Range("J18").Value = '0'
$ws.Range("K18").Value = '66.66'

# Row 19
$ws.Range("A19").Value = ' Sharjah'
$ws.Range("B19").Value = ' October 09 2020'
$ws.Range("C19").Value = 'Capitals won by 46 runs'
$ws.Range("D19").Value = 'Delhi Capitals'
$ws.Range("E19").Value = 'Rajasthan Royals'
$ws.Range("F19").Value = 'Axar Patel '
$ws.Range("G19").Value = '17'
$ws.Range("H19").Value = '8'
$ws.Range("I19").Value = '2'
$ws.Range("J19").Value = '1'
$ws.Range("K19").Value = '212.50'

